# Exportando os dados para o excel conforme resposta da aplicação web
# Preenche as colunas Status do Empréstimo (F), ID do Empréstimo (G) e APR (H)
# para cada pedido, com base no retorno da aplicação web de processamento de
# pedidos de empréstimo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pedido aprovado - maria@teste.com
$ws.Range("F2").Value = "APROVADO"
$ws.Range("G2").Value = "652347f1b75a6174ec3ac7da"
$ws.Range("H2").Value = 8

# Pedido aprovado - pedro@teste.com
$ws.Range("F3").Value = "APROVADO"
$ws.Range("G3").Value = "65234836b75a6174ec3ac7db"
$ws.Range("H3").Value = 4

# Pedido aprovado - daniela@teste.com
$ws.Range("F4").Value = "APROVADO"
$ws.Range("G4").Value = "6523487bb75a6174ec3ac7dc"
$ws.Range("H4").Value = 9

# Pedido aprovado - joao@teste.com
$ws.Range("F5").Value = "APROVADO"
$ws.Range("G5").Value = "652348c1b75a6174ec3ac7dd"
$ws.Range("H5").Value = 8

# Pedido não aprovado - sandra@teste.com
$ws.Range("F6").Value = "NÃO APROVADO!"

# Pedido não aprovado - marcelo@teste.com
$ws.Range("F7").Value = "NÃO APROVADO!"

# Seleciona o intervalo recém-preenchido, como ficaria após colar a resposta
$ws.Range("F2:H2").Select()
